# Auto-generated from the scraper diff: refresh market-price-derived
# columns (H:currentAveragePrice, I:currentAveragePriceNQ,
# J:currentAveragePriceHQ, K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ,
# N:LeveProfitHQ) on the affected rows of each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3018.4666
$ws.Range("I40").Value = 1796.8334
$ws.Range("J40").Value = 3832.889
$ws.Range("K40").Value = 1796.8334
$ws.Range("L40").Value = 3832.889
$ws.Range("M40").Value = -1621.8334
$ws.Range("N40").Value = -4182.889

$ws.Range("H69").Value = 500003520
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 500003520
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 1500010560
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -1500012308

$ws.Range("H72").Value = 500003520
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 500003520
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 4500031680
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -4500040416

$ws.Range("H86").Value = 4438.25
$ws.Range("J86").Value = 4376.75
$ws.Range("L86").Value = 4376.75
$ws.Range("N86").Value = -6622.75

$ws.Range("H89").Value = 4438.25
$ws.Range("J89").Value = 4376.75
$ws.Range("L89").Value = 21883.75
$ws.Range("N89").Value = -33115.75

$ws.Range("H94").Value = 2398.5
$ws.Range("I94").Value = 2398.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2398.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1947.5
$ws.Range("N94").ClearContents()

$ws.Range("H98").Value = 9330.869000000001
$ws.Range("I98").Value = 12815.6
$ws.Range("J98").Value = 2797
$ws.Range("K98").Value = 12815.6
$ws.Range("L98").Value = 2797
$ws.Range("M98").Value = -11317.6
$ws.Range("N98").Value = -5793

$ws.Range("H99").Value = 388.16666
$ws.Range("I99").Value = 388.16666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1164.49998
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 333.5000199999999
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 9330.869000000001
$ws.Range("I122").Value = 12815.6
$ws.Range("J122").Value = 2797
$ws.Range("K122").Value = 38446.8
$ws.Range("L122").Value = 8391
$ws.Range("M122").Value = -35996.8
$ws.Range("N122").Value = -13291

$ws.Range("H135").Value = 1498.1
$ws.Range("I135").Value = 1628.5555
$ws.Range("J135").Value = 324
$ws.Range("K135").Value = 14656.9995
$ws.Range("L135").Value = 2916
$ws.Range("M135").Value = -12121.9995
$ws.Range("N135").Value = -7986

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1677.5151
$ws.Range("I74").Value = 1763.1666
$ws.Range("K74").Value = 1763.1666
$ws.Range("M74").Value = -889.1666

$ws.Range("H77").Value = 1677.5151
$ws.Range("I77").Value = 1763.1666
$ws.Range("K77").Value = 8815.833000000001
$ws.Range("M77").Value = -4447.833000000001

$ws.Range("H88").Value = 2407.8572
$ws.Range("I88").Value = 1089
$ws.Range("J88").Value = 4166.3335
$ws.Range("K88").Value = 1089
$ws.Range("L88").Value = 4166.3335
$ws.Range("M88").Value = -683
$ws.Range("N88").Value = -4978.3335

$ws.Range("H91").Value = 2407.8572
$ws.Range("I91").Value = 1089
$ws.Range("J91").Value = 4166.3335
$ws.Range("K91").Value = 1089
$ws.Range("L91").Value = 4166.3335
$ws.Range("M91").Value = 315
$ws.Range("N91").Value = -6974.3335

$ws.Range("H96").Value = 63100
$ws.Range("J96").Value = 63100
$ws.Range("L96").Value = 63100
$ws.Range("N96").Value = -68592

$ws.Range("H132").Value = 2560.6562
$ws.Range("I132").Value = 1997.32
$ws.Range("K132").Value = 5991.96
$ws.Range("M132").Value = -3461.96

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1360.9474
$ws.Range("I20").Value = 1380.6364
$ws.Range("J20").Value = 1333.875
$ws.Range("K20").Value = 1380.6364
$ws.Range("L20").Value = 1333.875
$ws.Range("M20").Value = -1133.6364
$ws.Range("N20").Value = -1827.875

$ws.Range("H86").Value = 55558664
$ws.Range("I86").Value = 4333.3335
$ws.Range("K86").Value = 4333.3335
$ws.Range("M86").Value = -3210.3335

$ws.Range("H89").Value = 55558664
$ws.Range("I89").Value = 4333.3335
$ws.Range("K89").Value = 21666.6675
$ws.Range("M89").Value = -16050.6675

$ws.Range("H94").Value = 3607.077
$ws.Range("I94").Value = 3607.077
$ws.Range("K94").Value = 3607.077
$ws.Range("M94").Value = -3156.077

$ws.Range("H95").Value = 27170
$ws.Range("J95").Value = 27170
$ws.Range("L95").Value = 27170
$ws.Range("N95").Value = -32662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1214.931
$ws.Range("J22").Value = 1568.0769
$ws.Range("L22").Value = 1568.0769
$ws.Range("N22").Value = -2268.0769

$ws.Range("H28").Value = 74500
$ws.Range("J28").Value = 74500
$ws.Range("L28").Value = 74500
$ws.Range("N28").Value = -74990

$ws.Range("H38").Value = 6021
$ws.Range("J38").Value = 6021
$ws.Range("L38").Value = 6021
$ws.Range("N38").Value = -6775

$ws.Range("H46").Value = 6021
$ws.Range("J46").Value = 6021
$ws.Range("L46").Value = 6021
$ws.Range("N46").Value = -6443

$ws.Range("H58").Value = 2076.1667
$ws.Range("I58").Value = 2064.25
$ws.Range("J58").Value = 2100
$ws.Range("K58").Value = 2064.25
$ws.Range("L58").Value = 2100
$ws.Range("M58").Value = -1861.25
$ws.Range("N58").Value = -2506

$ws.Range("H132").Value = 2823.3157
$ws.Range("I132").Value = 3009
$ws.Range("K132").Value = 9027
$ws.Range("M132").Value = -6497

$ws.Range("H136").Value = 2076.1667
$ws.Range("I136").Value = 2064.25
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 6192.75
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -3642.75
$ws.Range("N136").Value = -11400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1365.475
$ws.Range("I34").Value = 180.66667
$ws.Range("J34").Value = 1461.5405
$ws.Range("K34").Value = 542.00001
$ws.Range("L34").Value = 4384.6215
$ws.Range("M34").Value = -458.00001
$ws.Range("N34").Value = -4552.6215

$ws.Range("H47").Value = 500
$ws.Range("J47").Value = 500
$ws.Range("L47").Value = 1500
$ws.Range("N47").Value = -2362

$ws.Range("H51").Value = 1400
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 3000
$ws.Range("M51").Value = -2540

$ws.Range("H92").Value = 3650
$ws.Range("I92").Value = 4250
$ws.Range("J92").Value = 2750
$ws.Range("K92").Value = 12750
$ws.Range("L92").Value = 8250
$ws.Range("M92").Value = -11502
$ws.Range("N92").Value = -10746

$ws.Range("H93").Value = 57013.5
$ws.Range("J93").Value = 57013.5
$ws.Range("L93").Value = 171040.5
$ws.Range("N93").Value = -174784.5

$ws.Range("H121").Value = 1103.6
$ws.Range("I121").Value = 815
$ws.Range("K121").Value = 2445
$ws.Range("M121").Value = -1135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 7949.875
$ws.Range("J92").Value = 7949.875
$ws.Range("L92").Value = 7949.875
$ws.Range("N92").Value = -11693.875

$ws.Range("H96").Value = 35999.8
$ws.Range("J96").Value = 35999.8
$ws.Range("L96").Value = 35999.8
$ws.Range("N96").Value = -41491.8

$ws.Range("H98").Value = 27571.5
$ws.Range("J98").Value = 27571.5
$ws.Range("L98").Value = 27571.5
$ws.Range("N98").Value = -33561.5

$ws.Range("H122").Value = 7011.4165
$ws.Range("I122").Value = 6785.227
$ws.Range("K122").Value = 20355.681
$ws.Range("M122").Value = -17905.681

$ws.Range("H132").Value = 3992.1052
$ws.Range("I132").Value = 3614.8125
$ws.Range("J132").Value = 6004.3335
$ws.Range("K132").Value = 10844.4375
$ws.Range("L132").Value = 18013.0005
$ws.Range("M132").Value = -8314.4375
$ws.Range("N132").Value = -23073.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1157.2106
$ws.Range("I16").Value = 616.82355
$ws.Range("J16").Value = 5750.5
$ws.Range("K16").Value = 616.82355
$ws.Range("L16").Value = 5750.5
$ws.Range("M16").Value = -446.82355
$ws.Range("N16").Value = -6090.5

$ws.Range("H46").Value = 1400.6
$ws.Range("J46").Value = 1400.6
$ws.Range("L46").Value = 1400.6
$ws.Range("N46").Value = -1776.6

$ws.Range("H122").Value = 4575.75
$ws.Range("I122").Value = 4236.2354
$ws.Range("K122").Value = 12708.7062
$ws.Range("M122").Value = -10258.7062

$ws.Range("H132").Value = 3798.1296
$ws.Range("I132").Value = 3929.0417
$ws.Range("K132").Value = 11787.1251
$ws.Range("M132").Value = -9257.125100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 33358332
$ws.Range("I5").Value = 44998
$ws.Range("K5").Value = 44998
$ws.Range("M5").Value = -44886

$ws.Range("H113").Value = 2441.5557
$ws.Range("I113").Value = 2162.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 6486.999899999999
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -4316.999899999999
$ws.Range("N113").Value = -13340
